# Fruta / hortaliza, semanal
# Update rows 2-9 and 11 (Fecha, Volumen, Precio minimo/maximo/promedio, Precio $/Kg)
# Row 10 is left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values (D, M, N, O, P, S) for each target row, taken from the
# re-ordered weekly data set.
$updates = @{
    2  = @{ D = 44320; M = 80;  N = 16000; O = 17000; P = 16500; S = 825  }
    3  = @{ D = 44708; M = 80;  N = 20000; O = 21000; P = 20500; S = 1025 }
    4  = @{ D = 44798; M = 80;  N = 21000; O = 22000; P = 21500; S = 1075 }
    5  = @{ D = 44893; M = 80;  N = 21000; O = 22000; P = 21625; S = 1081 }
    6  = @{ D = 44357; M = 100; N = 14000; O = 15000; P = 14500; S = 725  }
    7  = @{ D = 45092; M = 150; N = 24000; O = 25000; P = 24333; S = 1217 }
    8  = @{ D = 44792; M = 100; N = 21000; O = 22000; P = 21500; S = 1075 }
    9  = @{ D = 44533; M = 100; N = 16000; O = 17000; P = 16500; S = 825  }
    11 = @{ D = 44761; M = 100; N = 20000; O = 21000; P = 20500; S = 1025 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
